# Update output files (.json, .xlsx)
#
# 1) Replace every occurrence of the Polish text "nie dotyczy" ("not
#    applicable") with its English equivalent "N/A" across all worksheets.
# 2) Narrow a recurring set of "Ilość" (quantity) columns — one per
#    weekday block, repeated across all 5 day-groups on every sheet —
#    from their old widths (14 / 8 / 10 / 8) down to a uniform width of 7.

$wb = $excel.ActiveWorkbook

# Column indices (1-based) whose raw width must become 7, and the
# raw width they currently hold (kept here only for reference/debugging -
# we overwrite unconditionally since every sheet shares this layout).
$targetCols = @(4, 7, 13, 16, 24, 27, 33, 36, 44, 47, 53, 56, 64, 67, 73, 76, 84, 87, 93, 96)

# Excel's ColumnWidth (character-width units) vs. the raw OOXML <col width>
# attribute differ by a constant 0.83 offset in this workbook's font metric
# (width 7 raw <-> 6.17 ColumnWidth, width 15 raw <-> 14.17 ColumnWidth, ...).
$newColumnWidth = 6.17

foreach ($ws in $wb.Worksheets) {

    # --- narrow the target columns ---
    foreach ($colIdx in $targetCols) {
        $ws.Columns.Item($colIdx).ColumnWidth = $newColumnWidth
    }

    # --- replace "nie dotyczy" -> "N/A" in every used cell ---
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $rowOffset = $used.Row
    $colOffset = $used.Column

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($rowOffset + $r, $colOffset + $c)
            if ($cell.Value2 -eq "nie dotyczy") {
                $cell.Value2 = "N/A"
            }
        }
    }
}
